$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for row 3 / row 4 (Iteration 3)
$ws.Range("H3").Value = 0
$ws.Range("G4").Value = "1.2, 1.3, 2.1, 2.2"
$ws.Range("H4").Formula = "=3+5+8+3"
$ws.Range("I4").Formula = "=I3-H4"

# Update selection to match author's last cursor position
$ws.Range("I12").Select()
